$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of data to append at row 46 (text columns A-L, numeric columns M-N)
$row = 46

$textValues = @{
    "A" = "6132"
    "B" = "6/12/2025"
    "C" = "MERCEDES 370"
    "D" = "10"
    "E" = "807458394"
    "F" = "GESTION TELECENTRO"
    "G" = "Pendiente"
    "H" = "con fuente TLC "
    "I" = "1"
    "J" = "Cambio"
    "K" = "Fuente TLC"
    "L" = "Pasante"
}

foreach ($col in $textValues.Keys) {
    $cell = $ws.Range("$col$row")
    # Force text storage so numeric-looking strings (e.g. "6132", "10") are
    # not reinterpreted as numbers, matching the original inlineStr cells.
    $cell.NumberFormat = "@"
    $cell.Value = $textValues[$col]
    # Reset to the default "Normal" style so no stray style index is left on
    # the cell (keeps it consistent with the rest of the sheet's plain cells).
    $cell.Style = "Normal"
}

# Numeric coordinate columns
$ws.Range("M$row").Value = -58.484808
$ws.Range("N$row").Value = -34.630188
